$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258, shifting the existing rows 258:339 down to 259:340.
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the new data record.
$ws.Cells.Item(258, 1).Value2 = 9
$ws.Cells.Item(258, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(258, 3).Value2 = "Metropolitana"
$ws.Cells.Item(258, 4).Value2 = 44855
$ws.Cells.Item(258, 5).Value2 = 13
$ws.Cells.Item(258, 6).Value2 = 100112017
$ws.Cells.Item(258, 7).Value2 = "Apio"
$ws.Cells.Item(258, 8).Value2 = "Americana (o)"
$ws.Cells.Item(258, 9).Value2 = "Primera"
$ws.Cells.Item(258, 10).Value2 = 100
$ws.Cells.Item(258, 11).Value2 = 7000
$ws.Cells.Item(258, 12).Value2 = 8000
$ws.Cells.Item(258, 13).Value2 = 7600
$ws.Cells.Item(258, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(258, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(258, 16).Value2 = 1267
$ws.Cells.Item(258, 17).Value2 = 6
$ws.Cells.Item(258, 18).Value2 = "Hortaliza"

# Apply the same date style used by the other rows' date column (D).
$ws.Cells.Item(258, 4).NumberFormat = $ws.Cells.Item(259, 4).NumberFormat
